$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.450.48"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("D3").Value = "1.943.65"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'242.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "'0.617"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("D7").Value = "'58.40"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.47%  "
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("D10").Value = "'55.85"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "'0.0831"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.85%  "
$ws.Range("D12").Value = "'0.104"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").Value = "'21.54"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.04%  "
$ws.Range("D14").Value = "'0.822"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.44%  "
$ws.Range("D15").Value = "2.225.95"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Value = "'13.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.20%  "
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").Value = "1.946.76"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "36.305.47"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "'69.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("D21").Value = "0.0₃0860"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").Value = "'228.70"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").Value = "'5.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.74%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("D25").Value = "'2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.43%  "
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "'9.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.42%  "
$ws.Range("D28").Value = "'161.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.13%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").Value = "'0.130"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'19.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.11%  "
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("E32").Value = "  +1.68%  "
$ws.Range("D33").Value = "'4.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.81%  "
$ws.Range("D34").Value = "'0.0627"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("D35").Value = "'4.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").Value = "'6.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("E39").Value = "  -5.43%  "
$ws.Range("D40").Value = "'3.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'0.0980"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("E43").Value = "  -3.91%  "
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").Value = "'15.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "1.352.12"
$ws.Range("E46").Value = "  +1.63%  "
$ws.Range("E47").Value = "  -4.62%  "
$ws.Range("D48").Value = "'87.85"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.33%  "
$ws.Range("D49").Value = "'7.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.59%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("D51").Value = "'45.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.91%  "
